$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 217-218, shifting existing rows 217:242 down to 219:244.
$ws.Rows("217:218").Insert()

# New row 217 data
$ws.Cells.Item(217, 1).Value = 10
$ws.Cells.Item(217, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(217, 3).Value = "La Araucanía"
$ws.Cells.Item(217, 4).Value = 44449
$ws.Cells.Item(217, 5).Value = 9
$ws.Cells.Item(217, 6).Value = 100112023
$ws.Cells.Item(217, 7).Value = "Brócoli"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "Primera"
$ws.Cells.Item(217, 10).Value = 1900
$ws.Cells.Item(217, 11).Value = 800
$ws.Cells.Item(217, 12).Value = 900
$ws.Cells.Item(217, 13).Value = 866
$ws.Cells.Item(217, 14).Value = "`$/unidad"
$ws.Cells.Item(217, 15).Value = "Región Metropolitana"
$ws.Cells.Item(217, 16).Value = 866
$ws.Cells.Item(217, 17).Value = 1
$ws.Cells.Item(217, 18).Value = "Hortaliza"

# New row 218 data
$ws.Cells.Item(218, 1).Value = 10
$ws.Cells.Item(218, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(218, 3).Value = "La Araucanía"
$ws.Cells.Item(218, 4).Value = 44449
$ws.Cells.Item(218, 5).Value = 9
$ws.Cells.Item(218, 6).Value = 100112023
$ws.Cells.Item(218, 7).Value = "Brócoli"
$ws.Cells.Item(218, 8).Value = "Sin especificar"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 750
$ws.Cells.Item(218, 11).Value = 800
$ws.Cells.Item(218, 12).Value = 800
$ws.Cells.Item(218, 13).Value = 800
$ws.Cells.Item(218, 14).Value = "`$/unidad"
$ws.Cells.Item(218, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(218, 16).Value = 800
$ws.Cells.Item(218, 17).Value = 1
$ws.Cells.Item(218, 18).Value = "Hortaliza"
